# Update column F (dSF) values on Sheet1 to reflect the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 5
    4  = -3
    5  = -2
    7  = 4
    8  = -3
    9  = -1
    11 = 0
    12 = 0
    13 = 3
    14 = 4
    15 = 4
    16 = -2
    17 = 3
    18 = -4
    20 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
